# Add a new bullet point to the "Website Log" table, right after the
# "Added LinkedIn buttons for some pages" entry, documenting that image
# legends were added for the image carousel.

$d = $word.ActiveDocument

# Locate the paragraph that ends the "Added ... LinkedIn buttons for some
# pages" bullet. We search document-wide (rather than scoping a Range and
# reading Range.Paragraphs, which misbehaves on short/collapsed ranges) so
# the match is robust.
$anchorText = "buttons for some pages"
$anchorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$anchorText*") {
        $anchorPara = $p
        break
    }
}
if ($null -eq $anchorPara) {
    throw "Could not locate the anchor paragraph containing '$anchorText'"
}

# Insert a brand-new paragraph right after the anchor paragraph.
$anchorPara.Range.InsertParagraphAfter() | Out-Null
$newPara = $anchorPara.Next()
$newRange = $newPara.Range

# Give the new paragraph the same list formatting as its sibling bullets
# (ListParagraph style, single-level list numbering id 6, and the banding
# cnfStyle used throughout this table row) and set its text, all in one
# shot via a WordprocessingML fragment so the paragraph properties come out
# exactly right.
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:pPr>' +
         '<w:pStyle w:val="ListParagraph"/>' +
         '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr>' +
         '<w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>' +
       '</w:pPr>' +
       '<w:r><w:t>Added image legends for image carousel for projects</w:t></w:r>' +
       '</w:p>'

$newRange.InsertXML($xml) | Out-Null

Write-Host "Inserted new bullet: 'Added image legends for image carousel for projects'"
